$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 5000
$ws.Range("I64").Value = 5000
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 5000
$ws.Range("L64").Value = 0
$ws.Range("M64").Value = -4752
$ws.Range("N64").ClearContents()

$ws.Range("H67").Value = 5000
$ws.Range("I67").Value = 5000
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 5000
$ws.Range("L67").Value = 0
$ws.Range("M67").Value = -4142
$ws.Range("N67").ClearContents()

$ws.Range("H74").Value = 4037.5
$ws.Range("I74").Value = 4076.923
$ws.Range("J74").Value = 3866.6667
$ws.Range("K74").Value = 4076.923
$ws.Range("L74").Value = 3866.6667
$ws.Range("M74").Value = -3140.923
$ws.Range("N74").Value = -5738.6667

$ws.Range("H76").Value = 3011.842
$ws.Range("I76").Value = 3002.5
$ws.Range("J76").Value = 3061.6667
$ws.Range("K76").Value = 3002.5
$ws.Range("L76").Value = 3061.6667
$ws.Range("M76").Value = -2687.5
$ws.Range("N76").Value = -3691.6667

$ws.Range("H77").Value = 4037.5
$ws.Range("I77").Value = 4076.923
$ws.Range("J77").Value = 3866.6667
$ws.Range("K77").Value = 20384.615
$ws.Range("L77").Value = 19333.3335
$ws.Range("M77").Value = -15704.615
$ws.Range("N77").Value = -28693.3335

$ws.Range("H79").Value = 3011.842
$ws.Range("I79").Value = 3002.5
$ws.Range("J79").Value = 3061.6667
$ws.Range("K79").Value = 3002.5
$ws.Range("L79").Value = 3061.6667
$ws.Range("M79").Value = -1910.5
$ws.Range("N79").Value = -5245.6667

$ws.Range("H112").Value = 41668000
$ws.Range("J112").Value = 55557100
$ws.Range("L112").Value = 166671300
$ws.Range("N112").Value = -166673516

$ws.Range("H137").Value = 32261874
$ws.Range("I137").Value = 2043.2858
$ws.Range("J137").Value = 58828796
$ws.Range("K137").Value = 6129.857400000001
$ws.Range("L137").Value = 176486388
$ws.Range("M137").Value = -3579.857400000001
$ws.Range("N137").Value = -176491488

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 19040.828
$ws.Range("I32").Value = 14114.127
$ws.Range("K32").Value = 14114.127
$ws.Range("M32").Value = -13827.127

$ws.Range("H63").Value = 4902.9375
$ws.Range("I63").Value = 2021
$ws.Range("J63").Value = 6212.909
$ws.Range("K63").Value = 2021
$ws.Range("L63").Value = 6212.909
$ws.Range("M63").Value = -1335
$ws.Range("N63").Value = -7584.909

$ws.Range("H66").Value = 4902.9375
$ws.Range("I66").Value = 2021
$ws.Range("J66").Value = 6212.909
$ws.Range("K66").Value = 10105
$ws.Range("L66").Value = 31064.545
$ws.Range("M66").Value = -6673
$ws.Range("N66").Value = -37928.545

$ws.Range("H88").Value = 1753
$ws.Range("I88").Value = 1753
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 1753
$ws.Range("L88").Value = 0
$ws.Range("M88").Value = -1347
$ws.Range("N88").ClearContents()

$ws.Range("H91").Value = 1753
$ws.Range("I91").Value = 1753
$ws.Range("J91").Value = 0
$ws.Range("K91").Value = 1753
$ws.Range("L91").Value = 0
$ws.Range("M91").Value = -349
$ws.Range("N91").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 550
$ws.Range("I94").Value = 300
$ws.Range("J94").Value = 800
$ws.Range("K94").Value = 300
$ws.Range("L94").Value = 800
$ws.Range("M94").Value = 151
$ws.Range("N94").Value = -1702

$ws.Range("H105").Value = 1633.32
$ws.Range("I105").Value = 1512.1578
$ws.Range("J105").Value = 2017
$ws.Range("K105").Value = 1512.1578
$ws.Range("L105").Value = 2017
$ws.Range("M105").Value = 234.8422
$ws.Range("N105").Value = -5511

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6810
$ws.Range("I31").Value = 2885.5
$ws.Range("J31").Value = 7751.88
$ws.Range("K31").Value = 2885.5
$ws.Range("L31").Value = 7751.88
$ws.Range("M31").Value = -2590.5
$ws.Range("N31").Value = -8341.880000000001

$ws.Range("H34").Value = 6810
$ws.Range("I34").Value = 2885.5
$ws.Range("J34").Value = 7751.88
$ws.Range("K34").Value = 2885.5
$ws.Range("L34").Value = 7751.88
$ws.Range("M34").Value = -2683.5
$ws.Range("N34").Value = -8155.88

$ws.Range("H62").Value = 3785.7144
$ws.Range("I62").Value = 3300
$ws.Range("J62").Value = 5000
$ws.Range("K62").Value = 3300
$ws.Range("L62").Value = 5000
$ws.Range("M62").Value = -2676
$ws.Range("N62").Value = -6248

$ws.Range("H65").Value = 3785.7144
$ws.Range("I65").Value = 3300
$ws.Range("J65").Value = 5000
$ws.Range("K65").Value = 16500
$ws.Range("L65").Value = 25000
$ws.Range("M65").Value = -13380
$ws.Range("N65").Value = -31240

$ws.Range("H93").Value = 12631.4
$ws.Range("I93").Value = 10701.556
$ws.Range("K93").Value = 10701.556
$ws.Range("M93").Value = -8829.556

$ws.Range("H134").Value = 2354.8718
$ws.Range("I134").Value = 1912.9706
$ws.Range("K134").Value = 5738.9118
$ws.Range("M134").Value = -3203.9118

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 33772.363
$ws.Range("J37").Value = 33772.363
$ws.Range("L37").Value = 101317.089
$ws.Range("N37").Value = -101541.089

$ws.Range("H131").Value = 68634480
$ws.Range("J131").Value = 38466110
$ws.Range("L131").Value = 115398330
$ws.Range("N131").Value = -115408410

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4502.2
$ws.Range("I70").Value = 4277.7334
$ws.Range("K70").Value = 4277.7334
$ws.Range("M70").Value = -4007.7334

$ws.Range("H73").Value = 4502.2
$ws.Range("I73").Value = 4277.7334
$ws.Range("K73").Value = 4277.7334
$ws.Range("M73").Value = -3341.7334

$ws.Range("H80").Value = 87323.53999999999
$ws.Range("I80").Value = 3111.6667
$ws.Range("J80").Value = 159505.14
$ws.Range("K80").Value = 3111.6667
$ws.Range("L80").Value = 159505.14
$ws.Range("M80").Value = -2113.6667
$ws.Range("N80").Value = -161501.14

$ws.Range("H83").Value = 87323.53999999999
$ws.Range("I83").Value = 3111.6667
$ws.Range("J83").Value = 159505.14
$ws.Range("K83").Value = 15558.3335
$ws.Range("L83").Value = 797525.7000000001
$ws.Range("M83").Value = -10566.3335
$ws.Range("N83").Value = -807509.7000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 860.70966
$ws.Range("I46").Value = 889.4167
$ws.Range("K46").Value = 889.4167
$ws.Range("M46").Value = -701.4167

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 26803.273
$ws.Range("I136").Value = 30870.111
$ws.Range("K136").Value = 92610.333
$ws.Range("M136").Value = -90060.333
